# Repull data, push all data, mean calculation
# Update the "dSF" column (F) values for a set of rows to reflect the
# re-pulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -1
    9  = 0
    11 = -2
    12 = 1
    15 = 0
    17 = 1
    21 = 2
    23 = 2
    24 = 11
    30 = 1
    35 = -4
    36 = 1
    39 = 0
    49 = 3
    53 = -2
    56 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
